$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily observation was inserted above the current row 61, pushing
# every subsequent record (old rows 61..139) down by one row (to 62..140).
# Insert a blank row at 61 first, shifting everything below it down.
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with the new record's values.
$ws.Range("A61").Value = 11
$ws.Range("B61").Value = "Vega Monumental Concepción"
$ws.Range("C61").Value = "Bíobío"
$ws.Range("D61").Value = 45036
$ws.Range("E61").Value = 8
$ws.Range("F61").Value = "Fruta"
$ws.Range("G61").Value = 100103
$ws.Range("H61").Value = "Frutos de hueso (carozo)"
$ws.Range("I61").Value = 100103002
$ws.Range("J61").Value = "Ciruela"
$ws.Range("K61").Value = "Angeleno"
$ws.Range("L61").Value = "Primera"
$ws.Range("M61").Value = 70
$ws.Range("N61").Value = 13000
$ws.Range("O61").Value = 14000
$ws.Range("P61").Value = 13571
$ws.Range("Q61").Value = "$/bandeja 18 kilos granel"
$ws.Range("R61").Value = "Región de O'Higgins"
$ws.Range("S61").Value = 754
$ws.Range("T61").Value = 18
